# Updates cryptos list price (D) and volume (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, in diff order.
$updates = [ordered]@{
    'D2' = '42.276.92'
    'E2' = '  -2.93%  '
    'D3' = '2.217.43'
    'E3' = '  -2.10%  '
    'D4' = '1.01'
    'E4' = '  +0.26%  '
    'D5' = '107.52'
    'E5' = '  -10.66%  '
    'D6' = '296.01'
    'E6' = '  +11.63%  '
    'E7' = '  -2.96%  '
    'E8' = '  -0.17%  '
    'D9' = '0.597'
    'E9' = '  -3.70%  '
    'D10' = '43.53'
    'E10' = '  -8.53%  '
    'D11' = '0.0911'
    'E11' = '  -3.40%  '
    'D12' = '54.49'
    'E12' = '  +0.46%  '
    'D13' = '8.78'
    'E13' = '  -5.20%  '
    'D14' = '0.980'
    'E14' = '  +7.49%  '
    'E15' = '  -2.65%  '
    'E16' = '  -2.30%  '
    'D17' = '2.548.41'
    'E17' = '  -2.42%  '
    'D18' = '2.238.89'
    'E18' = '  -1.35%  '
    'D19' = '42.281.02'
    'E19' = '  -2.89%  '
    'E20' = '  +7.27%  '
    'E21' = '  -4.49%  '
    'E22' = '  +0.33%  '
    'D23' = '3.47'
    'E23' = '  +20.69%  '
    'D24' = '2.30'
    'E24' = '  -3.86%  '
    'D25' = '228.22'
    'E25' = '  -2.91%  '
    'E26' = '  -5.32%  '
    'E27' = '  -1.66%  '
    'D28' = '11.57'
    'E28' = '  -3.07%  '
    'E29' = '  -1.09%  '
    'D30' = '38.20'
    'E30' = '  -8.71%  '
    'E31' = '  -4.59%  '
    'D32' = '173.72'
    'E32' = '  +0.93%  '
    'E33' = '  -3.36%  '
    'E34' = '  -2.36%  '
    'D35' = '5.60'
    'E35' = '  -2.22%  '
    'D36' = '5.02'
    'E36' = '  +9.46%  '
    'D37' = '4.33'
    'E37' = '  +0.52%  '
    'E38' = '  -3.26%  '
    'D39' = '0.0367'
    'E39' = '  -2.27%  '
    'E40' = '  -4.02%  '
    'D41' = '2.42'
    'E41' = '  -4.76%  '
    'D42' = '71.86'
    'E42' = '  -3.37%  '
    'E43' = '  -2.31%  '
    'E44' = '  +0.09%  '
    'D45' = '12.57'
    'E45' = '  -9.28%  '
    'E46' = '  -4.89%  '
    'D47' = '5.40'
    'E47' = '  -6.17%  '
    'E48' = '  +4.58%  '
    'D49' = '103.16'
    'E49' = '  +1.85%  '
    'D50' = '8.40'
    'E50' = '  -1.44%  '
    'E51' = '  +4.05%  '
}

# Column-D numeric-looking strings (e.g. "1.01", "0.0911") must be forced to
# text so Excel keeps them as literal strings instead of coercing them to
# numbers (which would also silently drop meaningful trailing zeros).
$forceText = @(
    'D4'
    'D5'
    'D6'
    'D9'
    'D10'
    'D11'
    'D12'
    'D13'
    'D14'
    'D23'
    'D24'
    'D25'
    'D28'
    'D30'
    'D32'
    'D35'
    'D36'
    'D37'
    'D39'
    'D41'
    'D42'
    'D45'
    'D47'
    'D49'
    'D50'
)

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($forceText -contains $cellRef) {
        $range.NumberFormat = "@"
        $range.Value = $updates[$cellRef]
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cellRef]
    }
}
